# Append the latest Kaspa buy (run on 2026-02-20) as a new row at the
# bottom of the data table (row 28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Date column stores its values as plain text (e.g. "02/13/2026"),
# not real dates. Force text formatting before assigning the string so
# Excel doesn't auto-convert it into a date serial number, then restore
# the default "Normal" style so the new row matches the formatting of
# the existing data rows (no explicit style override).
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = "02/20/2026"
$ws.Range("A28").Style = "Normal"

$ws.Range("B28").Value = 1623.165999999997
$ws.Range("C28").Value = 0.0304959566674019
$ws.Range("D28").Value = 50
